$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 78206
$ws.Range("B3").Value = 95.28203104621129
$ws.Range("B4").Value = 8.550566236222211
$ws.Range("B5").Value = 36.09
